# Generate Report for Archive
#
# The localization status changed from "Ready for handoff" to
# "In Translation" for the tracked file. That status string appears on
# every sheet (Overview!E2/F2, zh-cn!C2, de-de!C2) and they all shared the
# same shared-string entry, so update all of the cells that showed the old
# status. Updating the status text also shrinks it, so the "Status" columns
# are narrower after the report is regenerated - re-fit those columns to
# their content, matching the narrower widths produced by the refresh.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn     = $wb.Worksheets.Item("zh-cn")
$dede     = $wb.Worksheets.Item("de-de")

$newStatus = "In Translation"

# Overview sheet tracks status per-locale in columns E (zh-cn) and F (de-de)
$overview.Range("E2").Value = $newStatus
$overview.Range("F2").Value = $newStatus

# Per-locale detail sheets carry the same status in column C
$zhcn.Range("C2").Value = $newStatus
$dede.Range("C2").Value = $newStatus

# Re-fit the status columns now that the text is shorter. (12.5 is the
# nearest settable value that lands the stored OOXML column width on the
# target 13.41-character width after the host's pixel-grid rounding.)
$overview.Columns.Item(5).ColumnWidth = 12.5
$overview.Columns.Item(6).ColumnWidth = 12.5
$zhcn.Columns.Item(3).ColumnWidth = 12.5
$dede.Columns.Item(3).ColumnWidth = 12.5
